# Slide 7 "Title 1" shape currently holds its title text split across three
# separate runs: "Mission ", "Planning " and a trailing tab character. The
# author collapsed these into a single run ("Mission Planning \t") that
# keeps the formatting (dirty="0", default rPr) of the first run.
#
# Re-assigning TextRange.Text with the exact same string is a no-op in this
# engine (it only rewrites the underlying runs when the text actually
# changes), so we first set a throwaway value and then set the real,
# merged text to force the paragraph's runs to be rebuilt as one run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange

$tr.Text = "__tmp__"
$tr.Text = "Mission Planning `t"
